{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n//\n// Change 1: \"Images are to be stored within a folder labeled \"images\"\"\n//           -> ... labeled \"assets\"\n// Change 2: Insert a new bulleted requirement, \"At least one file must be\n//           inside of a folder\", as the first item of the \"Other\n//           requirements:\" list (right after the \"Other requirements:\"\n//           heading paragraph and before the list's current first item).\n\nconst body = context.document.body;\n\n// --- Change 1: update the folder-name text -------------------------------\nconst oldText = 'Images are to be stored within a folder labeled \"images\"';\nconst newText = 'Images are to be stored within a folder labeled \"assets\"';\n\nconst imagesResults = body.search(oldText, { matchCase: true });\nimagesResults.load(\"items\");\nawait context.sync();\n\nif (imagesResults.items.length > 0) {\n  const imagesParagraph = imagesResults.items[0].paragraphs.getFirst();\n  // Replace the run text only; this keeps the paragraph's existing\n  // numbering/list formatting (numId 1004, ilvl 1) untouched.\n  imagesParagraph.insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// --- Change 2: add the new \"Other requirements\" bullet -------------------\nconst otherResults = body.search(\"Other requirements:\", { matchCase: true });\notherResults.load(\"items\");\nawait context.sync();\n\nif (otherResults.items.length > 0) {\n  const otherHeadingParagraph = otherResults.items[0].paragraphs.getFirst();\n\n  // The paragraph right after the heading is the current first bullet of\n  // the \"Other requirements\" list (numId 1005, ilvl 0). Inserting the new\n  // paragraph immediately before it makes the new paragraph inherit that\n  // same list formatting automatically.\n  const firstRequirementParagraph = otherHeadingParagraph.getNext();\n  firstRequirementParagraph.load(\"text\");\n  await context.sync();\n\n  firstRequirementParagraph.insertParagraph(\n    \"At least one file must be inside of a folder\",\n    Word.InsertLocation.before\n  );\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is the open document.\n#\n# Change 1: \"Images are to be stored within a folder labeled \"images\"\"\n#           -> ... labeled \"assets\"\n# Change 2: Insert a new bulleted requirement, \"At least one file must be\n#           inside of a folder\", as the first item of the \"Other\n#           requirements:\" list (right after the \"Other requirements:\"\n#           heading paragraph and before the list's current first item).\n\n$d = $word.ActiveDocument\n\n# --- Change 1: update the folder-name text --------------------------------\n$oldText = 'Images are to be stored within a folder labeled \"images\"'\n$newText = 'Images are to be stored within a folder labeled \"assets\"'\n\nforeach ($para in $d.Paragraphs) {\n    if ($para.Range.Text -eq ($oldText + [char]13)) {\n        # Assigning .Text (without the trailing paragraph mark) replaces just\n        # the paragraph's run content, leaving its numbering/list formatting\n        # (numId 1004, ilvl 1) untouched.\n        $para.Range.Text = $newText\n        break\n    }\n}\n\n# --- Change 2: add the new \"Other requirements\" bullet --------------------\n$headingPara = $null\nforeach ($para in $d.Paragraphs) {\n    if ($para.Range.Text -eq (\"Other requirements:\" + [char]13)) {\n        $headingPara = $para\n        break\n    }\n}\n\nif ($headingPara -ne $null) {\n    # The paragraph right after the heading is the current first bullet of\n    # the \"Other requirements\" list (numId 1005, ilvl 0). Inserting the new\n    # paragraph immediately before it makes the new paragraph inherit that\n    # same list formatting automatically.\n    $firstRequirementPara = $headingPara.Next()\n    $insertionRange = $firstRequirementPara.Range\n    $insertionRange.InsertParagraphBefore()\n\n    $newPara = $headingPara.Next()\n    $newPara.Range.Text = \"At least one file must be inside of a folder\"\n}\n"}
